# "small problems part 2 July 18"
# Lesson-5-Advanced-JavaScript-Collections / Working with Callback Functions.xlsx
#
# The table on Sheet1 documents the "outer/inner/third" nested-map walkthrough.
# This edit:
#   1. Tightens several wrapped-text row heights (they were taller than needed).
#   2. Fixes row 35 (previously a duplicate of the "callback execution" row) to
#      correctly read "Outer callback execution".
#   3. Fills out row 36 (which only had column A filled in) and appends five
#      new rows (37-41) describing the inner/third callback executions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Row height corrections (wrapped text needed less vertical space) ---
$ws.Rows.Item(5).RowHeight  = 68
$ws.Rows.Item(6).RowHeight  = 85
$ws.Rows.Item(12).RowHeight = 68
$ws.Rows.Item(17).RowHeight = 51
$ws.Rows.Item(18).RowHeight = 51
$ws.Rows.Item(19).RowHeight = 68
$ws.Rows.Item(24).RowHeight = 51
$ws.Rows.Item(25).RowHeight = 51

# --- 2. Fix row 35 content ---
$ws.Range("A35").Value = "Outer callback execution"
$ws.Range("B35").Value = "[[1, 2], [3, 4]] and [5, 6]"

# --- 3. Complete row 36 and add new rows 36-41 ---
$ws.Range("B36").Value = "[[1, 2], [3, 4]] and [5, 6]"
$ws.Range("C36").Value = "None"
$ws.Range("D36").Value = "New Array"
$ws.Range("E36").Value = "Yes, explicitly returned by callback"

$ws.Rows.Item(37).RowHeight = 34
$ws.Range("A37").Value = "inner callback execution"
$ws.Range("B37").Value = "[1, 2] and [3, 4] and 5 and 6"
$ws.Range("C37").Value = "None"
$ws.Range("D37").Value = "number or new array"
$ws.Range("E37").Value = "Yes, used by map for transformation"

$ws.Rows.Item(38).RowHeight = 34
$ws.Range("A38").Value = "(===)"
$ws.Range("B38").Value = "[1, 2] and [3, 4] and 5 and 6"
$ws.Range("C38").Value = "None"
$ws.Range("D38").Value = "Boolean"
$ws.Range("E38").Value = "Yes, used by callback execution"

$ws.Rows.Item(39).RowHeight = 34
$ws.Range("A39").Value = "(+)"
$ws.Range("B39").Value = "5 and 6"
$ws.Range("C39").Value = "None"
$ws.Range("D39").Value = "Number"
$ws.Range("E39").Value = "Yes, used by callback execution"

$ws.Rows.Item(40).RowHeight = 34
$ws.Range("A40").Value = "third methd call (map)"
$ws.Range("B40").Value = "[1, 2] and [3, 4]"
$ws.Range("C40").Value = "None"
$ws.Range("D40").Value = "New Array"
$ws.Range("E40").Value = "Yes, used by inner callback execution"

$ws.Rows.Item(41).RowHeight = 34
$ws.Range("A41").Value = "third callback execution"
$ws.Range("B41").Value = "[1, 2] and [3, 4]"
$ws.Range("C41").Value = "None"
$ws.Range("D41").Value = "New arrays"
$ws.Range("E41").Value = "Yes, used by map for transformation"

# --- 4. Update view: scroll position + active selection ---
$ws.Range("E42").Select()
